# This script reflows three single-run paragraphs into multi-line runs
# by inserting manual line breaks (w:br) at the points shown in the target diff,
# using Find/Replace with the "^l" (manual line break) replacement code.
$d = $word.ActiveDocument

# Paragraph: Portuguese "Programa" bullet list -> split into lines with manual breaks
$found = $d.Content.Find.Execute(
    "- Características das águas de interesse para o tratamento: características físicas, químicas e bacteriológicas; - Padrão de Potabilidade; - Tecnologias de Tratamento de Água;- Unidades Constituintes de um Sistema de Abastecimento de Água;- Captação de Águas Subterrâneas e Captação de Águas Superficiais - Gradeamento, remoção de areia, casa de bombas; - Reservação; - Sistema de Tratamento de Água de Ciclo Completo; - Coagulação-floculação e Mistura Rápida; - Decantação: decantação convencional e de alta taxa e sistema de remoção de lodo;- Mecanismos da filtração, materiais filtrantes e fundos de filtros, hidráulica da filtração, filtração com taxa constante e taxa declinante, - Desinfecção: principais desinfetantes, cloração e cloro-amoniação, pré e pós-cloração, - Tratamento de resíduos gerados nas ETAs e reuso de água recuperada - A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Características das águas de interesse para o tratamento: características físicas, químicas e bacteriológicas; ^l- Padrão de Potabilidade; ^l- Tecnologias de Tratamento de Água;^l- Unidades Constituintes de um Sistema de Abastecimento de Água;^l- Captação de Águas Subterrâneas e Captação de Águas Superficiais ^l- Gradeamento, remoção de areia, casa de bombas; ^l- Reservação; ^l- Sistema de Tratamento de Água de Ciclo Completo; ^l- Coagulação-floculação e Mistura Rápida; ^l- Decantação: decantação convencional e de alta taxa e sistema de remoção de lodo;^l- Mecanismos da filtração, materiais filtrantes e fundos de filtros, hidráulica da filtração, filtração com taxa constante e taxa declinante, ^l- Desinfecção: principais desinfetantes, cloração e cloro-amoniação, pré e pós-cloração, ^l- Tratamento de resíduos gerados nas ETAs e reuso de água recuperada ^l- A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina.", 2)
if (-not $found) { throw "Find/Replace 0 did not match expected text" }
Write-Output "Replace 0: $found"

# Paragraph: English "Programa" bullet list (italic) -> split into lines with manual breaks
$found = $d.Content.Find.Execute(
    "- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;- Potability Standard;- Water Treatment Technologies;- Constituent Units of a Water Supply System;- Groundwater Catchment and Surface Water Catchment- Railing, sand removal, pump room;- Reservation;- Full Cycle Water Treatment System;- Coagulation-flocculation and Rapid Mixing;- Decantation: conventional and high rate decantation and sludge removal system;- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,- Treatment of waste generated at stations and reuse of recovered water- The discipline may have didactic trips to complement the content of the discipline",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;^l- Potability Standard;^l- Water Treatment Technologies;^l- Constituent Units of a Water Supply System;^l- Groundwater Catchment and Surface Water Catchment^l- Railing, sand removal, pump room;^l- Reservation;^l- Full Cycle Water Treatment System;^l- Coagulation-flocculation and Rapid Mixing;^l- Decantation: conventional and high rate decantation and sludge removal system;^l- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,^l- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,^l- Treatment of waste generated at stations and reuse of recovered water^l- The discipline may have didactic trips to complement the content of the discipline", 2)
if (-not $found) { throw "Find/Replace 1 did not match expected text" }
Write-Output "Replace 1: $found"

# Paragraph: Bibliografia citation block -> split into lines with manual breaks
$found = $d.Content.Find.Execute(
    "PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Tratamento de Água deAbastecimento por Filtração em Múltiplas Etapas. ASSOCIAÇÃO BRASILEIRA DEENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999. (Coordenação: Luiz DiBernardo). PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Noções Gerais de Tratamento eDisposição Final de Lodos e Estações de Tratamento de Água. ASSOCIAÇÃOBRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 200(Coordenação: Marco A.P. Reali).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "PROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Tratamento de Água de^lAbastecimento por Filtração em Múltiplas Etapas. ASSOCIAÇÃO BRASILEIRA DE^lENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999. (Coordenação: Luiz Di^lBernardo). ^lPROGRAMA DE PESQUISA EM SANEAMENTO BÁSICO. Noções Gerais de Tratamento e^lDisposição Final de Lodos e Estações de Tratamento de Água. ASSOCIAÇÃO^lBRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 200^l(Coordenação: Marco A.P. Reali).", 2)
if (-not $found) { throw "Find/Replace 2 did not match expected text" }
Write-Output "Replace 2: $found"

